# NIT-9012107141.xlsx - Estado de Cuenta update
# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The previous account-statement rows (7 rows / 4 workers / 6 periods) are
# replaced by a new set of rows (19 rows / 5 workers / 17 periods), and the
# summary header totals are refreshed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Make room for the new data. The table currently has 7 data rows
#    (16-22); the new data needs 19 rows (16-34), so insert 12 blank rows
#    right before the last existing data row (22) - this pushes the old
#    row 22 (and the footer block below it) down to make space, while
#    rows 16-21 stay put.
# ---------------------------------------------------------------------------
$ws.Rows(22).Resize(12).Insert()

# Copy the normal (non-bottom-border) data-row formatting from row 21 onto
# the freshly inserted rows 22:33 so they look like the rest of the table
# (the old row 22 - now shifted to row 34 - already keeps the special
# bottom-border formatting for the last row of the table).
$ws.Range("B21:J21").Copy()
$ws.Range("B22:J33").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Write the new data set (rows 16-34).
# ---------------------------------------------------------------------------
$data = @(
    @("CC", "1143367398", "JANER LATORRE SALCEDO",             "2106", 20593, 908526),
    @("CC", "1047448271", "JULIETH PAOLA JIMENEZ MALDONADO",   "1903",  7729, 828116),
    @("CC", "1047448271", "JULIETH PAOLA JIMENEZ MALDONADO",   "1902", 33125, 828116),
    @("CC", "1001968718", "JOSE DAVID MESTRE TORRES",          "2408", 46400, 1160000),
    @("CC", "1001968718", "JOSE DAVID MESTRE TORRES",          "2407", 46400, 1160000),
    @("CC", "1001968718", "JOSE DAVID MESTRE TORRES",          "2406", 46400, 1160000),
    @("CC", "1001968718", "JOSE DAVID MESTRE TORRES",          "2405", 46400, 1160000),
    @("CC", "1001968718", "JOSE DAVID MESTRE TORRES",          "2404", 46400, 1160000),
    @("CC", "1001968718", "JOSE DAVID MESTRE TORRES",          "2403", 46400, 1160000),
    @("CC", "1001968718", "JOSE DAVID MESTRE TORRES",          "2402", 46400, 1160000),
    @("CC", "1001968718", "JOSE DAVID MESTRE TORRES",          "2401", 46400, 1160000),
    @("CC", "1001968718", "JOSE DAVID MESTRE TORRES",          "2312", 46400, 1160000),
    @("CC", "1001968718", "JOSE DAVID MESTRE TORRES",          "2311", 46400, 1160000),
    @("CC", "1001968718", "JOSE DAVID MESTRE TORRES",          "2310", 46400, 1160000),
    @("CC", "1001968718", "JOSE DAVID MESTRE TORRES",          "2309", 46400, 1160000),
    @("CC", "1017169516", "ANDRES FELIPE ARISTIZABAL GIRALDO", "2107", 36341, 908526),
    @("CC", "1017169516", "ANDRES FELIPE ARISTIZABAL GIRALDO", "2106", 36341, 908526),
    @("CC", "1001976349", "YAN CARLOS POLO CORPAS",            "2409", 52000, 1300000),
    @("CC", "1001976349", "YAN CARLOS POLO CORPAS",            "2408", 17333, 1300000)
)

$row = 16
foreach ($rec in $data) {
    $ws.Cells.Item($row, 2).Value = $rec[0]   # B - Tipo Doc Trabajador
    $ws.Cells.Item($row, 3).Value = $rec[1]   # C - N Doc Trabajador
    $ws.Cells.Item($row, 4).Value = $rec[2]   # D - Nombre Trabajador
    $ws.Cells.Item($row, 5).Value = $rec[3]   # E - Periodo Mora
    $ws.Cells.Item($row, 6).Value = $rec[4]   # F - Valor Mora
    $ws.Cells.Item($row, 7).Value = $rec[5]   # G - Salario Basico
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# 3. Refresh the summary header: total overdue value, worker count and
#    period count.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 760262
$ws.Range("C13").Value = 5
$ws.Range("F13").Value = 17

# ---------------------------------------------------------------------------
# 4. Re-fit the data columns now that longer names/numbers are present.
# ---------------------------------------------------------------------------
$ws.Columns("B:J").AutoFit()
